$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 102.4929963333333
$ws.Range("H2").Value = 307.478989
$ws.Range("I2").Value = 0.2065071987599813
$ws.Range("J2").Value = 0.2065071987599814
$ws.Range("M2").Value = 83.91225566666667
$ws.Range("N2").Value = 251.736767
$ws.Range("O2").Value = 0.9556261553553385
$ws.Range("P2").Value = 0.9556261553553385
$ws.Range("Q2").Value = 8600.418512365397
$ws.Range("R2").Value = 77403.76661128856
$ws.Range("S2").Value = 0.1973436804042017
$ws.Range("T2").Value = 0.1973436804042017

# Row 3
$ws.Range("G3").Value = 102.4929963333333
$ws.Range("H3").Value = 307.478989
$ws.Range("I3").Value = 0.2065071987599813
$ws.Range("J3").Value = 0.2065071987599814
$ws.Range("O3").Value = 0.00439999103960854
$ws.Range("P3").Value = 0.00439999103960854
$ws.Range("Q3").Value = 39.59892074868979
$ws.Range("R3").Value = 356.390286738208
$ws.Range("S3").Value = 0.0009086298241585777
$ws.Range("T3").Value = 0.0009086298241585779

# Row 4
$ws.Range("G4").Value = 102.4929963333333
$ws.Range("H4").Value = 307.478989
$ws.Range("I4").Value = 0.2065071987599813
$ws.Range("J4").Value = 0.2065071987599814
$ws.Range("M4").Value = 3.510050666666667
$ws.Range("N4").Value = 10.530152
$ws.Range("O4").Value = 0.03997385360505296
$ws.Range("P4").Value = 0.03997385360505297
$ws.Range("Q4").Value = 359.7556101084809
$ws.Range("R4").Value = 3237.800490976329
$ws.Range("S4").Value = 0.008254888531621069
$ws.Range("T4").Value = 0.00825488853162107

# Row 5
$ws.Range("I5").Value = 0.581825957350084
$ws.Range("J5").Value = 0.5818259573500841
$ws.Range("M5").Value = 83.91225566666667
$ws.Range("N5").Value = 251.736767
$ws.Range("O5").Value = 0.9556261553553385
$ws.Range("P5").Value = 0.9556261553553385
$ws.Range("Q5").Value = 24231.34285204438
$ws.Range("R5").Value = 218082.0856683995
$ws.Range("S5").Value = 0.5560081027084
$ws.Range("T5").Value = 0.5560081027084001

# Row 6
$ws.Range("I6").Value = 0.581825957350084
$ws.Range("J6").Value = 0.5818259573500841
$ws.Range("O6").Value = 0.00439999103960854
$ws.Range("P6").Value = 0.00439999103960854
$ws.Range("S6").Value = 0.00256002899895203
$ws.Range("T6").Value = 0.002560028998952031

# Row 7
$ws.Range("I7").Value = 0.581825957350084
$ws.Range("J7").Value = 0.5818259573500841
$ws.Range("M7").Value = 3.510050666666667
$ws.Range("N7").Value = 10.530152
$ws.Range("O7").Value = 0.03997385360505296
$ws.Range("P7").Value = 0.03997385360505297
$ws.Range("Q7").Value = 1013.597363773806
$ws.Range("R7").Value = 9122.376273964255
$ws.Range("S7").Value = 0.02325782564273205
$ws.Range("T7").Value = 0.02325782564273205

# Row 8
$ws.Range("G8").Value = 105.053815
$ws.Range("H8").Value = 315.161445
$ws.Range("I8").Value = 0.2116668438899346
$ws.Range("J8").Value = 0.2116668438899346
$ws.Range("M8").Value = 83.91225566666667
$ws.Range("N8").Value = 251.736767
$ws.Range("O8").Value = 0.9556261553553385
$ws.Range("P8").Value = 0.9556261553553385
$ws.Range("Q8").Value = 8815.302583038701
$ws.Range("R8").Value = 79337.72324734832
$ws.Range("S8").Value = 0.2022743722427368
$ws.Range("T8").Value = 0.2022743722427369

# Row 9
$ws.Range("G9").Value = 105.053815
$ws.Range("H9").Value = 315.161445
$ws.Range("I9").Value = 0.2116668438899346
$ws.Range("J9").Value = 0.2116668438899346
$ws.Range("O9").Value = 0.00439999103960854
$ws.Range("P9").Value = 0.00439999103960854
$ws.Range("Q9").Value = 40.58831181989334
$ws.Range("R9").Value = 365.29480637904
$ws.Range("S9").Value = 0.000931332216497932
$ws.Range("T9").Value = 0.0009313322164979321

# Row 10
$ws.Range("G10").Value = 105.053815
$ws.Range("H10").Value = 315.161445
$ws.Range("I10").Value = 0.2116668438899346
$ws.Range("J10").Value = 0.2116668438899346
$ws.Range("M10").Value = 3.510050666666667
$ws.Range("N10").Value = 10.530152
$ws.Range("O10").Value = 0.03997385360505296
$ws.Range("P10").Value = 0.03997385360505297
$ws.Range("Q10").Value = 368.7442133766267
$ws.Range("R10").Value = 3318.69792038964
$ws.Range("S10").Value = 0.008461139430699845
$ws.Range("T10").Value = 0.008461139430699849
